$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column J ("Candida genome database description and notes") with the
# expanded CGD descriptions (header + one entry per data row, folding the old
# free-floating row 15 note into row 14's cell). Cells are touched in the same
# order the author filled them in (J2/J3 before circling back to the J1
# header) so newly-introduced shared strings land in the same pool order.
$ws.Range("J2").Value = "Stationary phase protein; vitamin B synthesis; induced byyeast-hypha switch, 3-AT or in azole-resistant strain overexpressing MDR1; soluble in hyphae; regulated by Gcn4, macrophage; Spider biofilm induced; rat catheter biofilm repressed. Notes: pyrioxidine is the 4-methol version of vitamin B6."
$ws.Range("J3").Value = "Protein with a predicted role in pyridoxine metabolism; stationary phase protein; regulated by Tup1, Efg1; Spider biofilm induced. Notes: Pyridoxal 5'-phosphate (PLP) is the active form of vitamin B6."
$ws.Range("J1").Value = "Candida genome database description and notes "
$ws.Range("J4").Value = "Putative trifunctional enzyme of thiamine biosynthesis, degradation and salvage; Spider biofilm induced. Notes: THI20 adds phosphate groups to pyrimidines and breaks down thiamine."
$ws.Range("J5").Value = "Putative farnesyl pyrophosphate synthetase involved in isoprenoid and sterol biosynthesis, based on similarity to S. cerevisiae Erg20p; likely to be essential for growth, based on an insertional mutagenesis strategy. Notes: This enzyme is important to making thiamine precurors."
$ws.Range("J6").Value = "Putative thiamin-phosphate pyrophosphorylase, hydroxyethylthiazole kinase; fungal-specific; Spider biofilm induced. Notes: This enzyme makes thiamine diphosphates."
$ws.Range("J7").Value = "None: uncharacterized."
$ws.Range("J8").Value = "Thiamine biosynthetic enzyme precursor; repressed during the mating process; stationary phase enriched protein; Spider biofilm induced. Notes: Thiamine is a thiazole linked to a pyridimine."
$ws.Range("J9").Value = "Subunit of the Dam1 (DASH) complex, which acts in chromosome segregation by coupling kinetochores to spindle microtubules . Notes: Miotic spindle protein"
$ws.Range("J10").Value = "None: uncharacterized."
$ws.Range("J11").Value = "None: uncharacterized."
$ws.Range("J12").Value = "Thiamin pyrimidine synthase; synthesis of the thiamine precursor hydroxymethylpyrimidine phosphate; single-turnover enzyme that provides histidine for HMP-P formation; induced by nitric oxide independent of Yhb1; Spider biofilm induced. NotesL This enzyme makes HMP-P, which is an important precursor of TPP synthesis, which is a form of vitamine B."
$ws.Range("J13").Value = "None: uncharacterized."
$ws.Range("J14").Value = "Putative transporter; more similar to S. cerevisiae Tpn1, which is a vitamin B6 transporter, than to purine-cytosine permeases; transcription is regulated by Nrg1; Spider biofilm induced. Notes: possible role in Vitamin B transport."

# Row 2 grew a touch taller to fit the new wrapped text.
$ws.Rows.Item(2).RowHeight = 21

# Drop the now-empty trailing row (its only cell, J15, is folded into J14 above).
$ws.Range("J15").ClearContents() | Out-Null

# Leave the selection where the last edit landed.
$ws.Range("D14").Select() | Out-Null
